$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d2Style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.339.94"
$ws.Range("D2").Style = $d2Style
$ws.Range("E2").Value = "  +5.99%  "

$d3Style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.039.65"
$ws.Range("D3").Style = $d3Style
$ws.Range("E3").Value = "  +8.03%  "

$d5Style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.72"
$ws.Range("D5").Style = $d5Style
$ws.Range("E5").Value = "  +2.98%  "

$d6Style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.691"
$ws.Range("D6").Style = $d6Style
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("E7").Value = "  -0.02%  "

$d8Style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.65"
$ws.Range("D8").Style = $d8Style
$ws.Range("E8").Value = "  +9.08%  "

$d9Style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("D9").Style = $d9Style
$ws.Range("E9").Value = "  +7.86%  "

$d10Style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.25"
$ws.Range("D10").Style = $d10Style
$ws.Range("E10").Value = "  +4.76%  "

$d11Style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0761"
$ws.Range("D11").Style = $d11Style
$ws.Range("E11").Value = "  +2.66%  "

$d12Style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.100"
$ws.Range("D12").Style = $d12Style
$ws.Range("E12").Value = "  +2.71%  "

$d13Style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.39"
$ws.Range("D13").Style = $d13Style
$ws.Range("E13").Value = "  +11.94%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$d14Style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.333.08"
$ws.Range("D14").Style = $d14Style
$ws.Range("E14").Value = "  +8.00%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$d15Style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.835"
$ws.Range("D15").Style = $d15Style
$ws.Range("E15").Value = "  +6.68%  "

$d16Style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("D16").Style = $d16Style
$ws.Range("E16").Value = "  +4.65%  "

$d17Style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.037.53"
$ws.Range("D17").Style = $d17Style
$ws.Range("E17").Value = "  +8.40%  "

$d18Style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.432.85"
$ws.Range("D18").Style = $d18Style
$ws.Range("E18").Value = "  +6.17%  "

$d19Style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.15"
$ws.Range("D19").Style = $d19Style
$ws.Range("E19").Value = "  +2.66%  "

$d20Style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0860"
$ws.Range("D20").Style = $d20Style
$ws.Range("E20").Value = "  +4.51%  "

$ws.Range("E21").Value = "  +7.87%  "

$d22Style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "252.97"
$ws.Range("D22").Style = $d22Style
$ws.Range("E22").Value = "  +3.81%  "

$d23Style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.27"
$ws.Range("D23").Style = $d23Style
$ws.Range("E23").Value = "  +2.20%  "

$d24Style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = $d24Style
$ws.Range("E24").Value = "  +0.09%  "

$d25Style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = $d25Style
$ws.Range("E25").Value = "  -4.57%  "

$d26Style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.39"
$ws.Range("D26").Style = $d26Style
$ws.Range("E26").Value = "  +1.89%  "

$d27Style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.10"
$ws.Range("D27").Style = $d27Style
$ws.Range("E27").Value = "  -1.82%  "

$d28Style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.69"
$ws.Range("D28").Style = $d28Style
$ws.Range("E28").Value = "  +13.43%  "

$d29Style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.93"
$ws.Range("D29").Style = $d29Style
$ws.Range("E29").Value = "  +5.39%  "

$ws.Range("E30").Value = "  +2.03%  "

$d31Style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.61"
$ws.Range("D31").Style = $d31Style
$ws.Range("E31").Value = "  +70.69%  "

$d32Style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.57"
$ws.Range("D32").Style = $d32Style
$ws.Range("E32").Value = "  +6.47%  "

$d33Style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0615"
$ws.Range("D33").Style = $d33Style
$ws.Range("E33").Value = "  +3.65%  "

$d34Style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.33"
$ws.Range("D34").Style = $d34Style
$ws.Range("E34").Value = "  +4.04%  "

$d35Style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0890"
$ws.Range("D35").Style = $d35Style
$ws.Range("E35").Value = "  +24.87%  "

$ws.Range("E36").Value = "  +0.01%  "

$d37Style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").Style = $d37Style
$ws.Range("E37").Value = "  +1.41%  "

$d38Style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.25"
$ws.Range("D38").Style = $d38Style
$ws.Range("E38").Value = "  +16.42%  "

$d39Style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.906"
$ws.Range("D39").Style = $d39Style
$ws.Range("E39").Value = "  +6.87%  "

$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("E41").Value = "  +4.16%  "

$d42Style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.32"
$ws.Range("D42").Style = $d42Style
$ws.Range("E42").Value = "  +3.53%  "

$d43Style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.22"
$ws.Range("D43").Style = $d43Style
$ws.Range("E43").Value = "  +1.21%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$d44Style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.14"
$ws.Range("D44").Style = $d44Style
$ws.Range("E44").Value = "  +6.66%  "

$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$d45Style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.88"
$ws.Range("D45").Style = $d45Style
$ws.Range("E45").Value = "  +19.50%  "

$d46Style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.374.66"
$ws.Range("D46").Style = $d46Style
$ws.Range("E46").Value = "  +3.72%  "

$d48Style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("D48").Style = $d48Style
$ws.Range("E48").Value = "  +2.03%  "

$ws.Range("E49").Value = "  +4.50%  "

$d50Style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.223.21"
$ws.Range("D50").Style = $d50Style
$ws.Range("E50").Value = "  +7.91%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$d51Style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.93"
$ws.Range("D51").Style = $d51Style
$ws.Range("E51").Value = "  +19.03%  "
